# Applies the OOXML diff: rewrites part of the Model paragraph's
# rationale text (inheritance -> association / singleton motivation),
# and splits a run mid-word to relocate the "_GoBack" bookmark.

$d = $word.ActiveDocument

function Force-Split($s, $e) {
    # Force the [s,e) text to live in run(s) distinct from whatever
    # precedes/follows it, by toggling a character-formatting property
    # on it and back off. No-op on the text; only affects how the
    # engine segments runs when it re-serializes.
    if ($e -gt $s) {
        $rr = $d.Range($s, $e)
        $rr.Bold = 1
        $rr.Bold = 0
    }
}

function Force-SplitAt($pos) {
    if ($pos -gt 0) {
        $rr = $d.Range($pos - 1, $pos)
        $rr.Bold = 1
        $rr.Bold = 0
    }
}

# ---------------------------------------------------------------------
# Edit 1: "You may ask why we have chosen inheritance over composition. "
#      -> "You may ask why we have chosen " / "association" / " over composition. "
# ---------------------------------------------------------------------
$full = $d.Content.Text
$oldText = "You may ask why we have chosen inheritance over composition. "
$idx = $full.IndexOf($oldText)
$newText = "You may ask why we have chosen association over composition. "
$r = $d.Range($idx, $idx + $oldText.Length)
$r.Text = $newText

$p1 = "You may ask why we have chosen "
$p2 = "association"
$p3 = " over composition. "
$s1 = $idx; $e1 = $s1 + $p1.Length
$s2 = $e1;  $e2 = $s2 + $p2.Length
$s3 = $e2;  $e3 = $s3 + $p3.Length

Force-Split $s1 $e1
Force-Split $s2 $e2
Force-Split $s3 $e3

# ---------------------------------------------------------------------
# Edit 2: " when " / "a " / "change is made on an object in one class,
#          it is expected to change in all of them"
#      -> " " / "we wanted to have only one instance of every object"
# ---------------------------------------------------------------------
$full = $d.Content.Text
$anchor = "is that when a change is made on an object in one class, it is expected to change in all of them."
$anchorIdx = $full.IndexOf($anchor)
$prefixLen = "is that".Length
$oldPart = " when a change is made on an object in one class, it is expected to change in all of them"
$start = $anchorIdx + $prefixLen
$newPart = " we wanted to have only one instance of every object"
$r = $d.Range($start, $start + $oldPart.Length)
$r.Text = $newPart

$p1 = " "
$p2 = "we wanted to have only one instance of every object"
$s1 = $start; $e1 = $s1 + $p1.Length
$s2 = $e1;    $e2 = $s2 + $p2.Length

Force-Split $s1 $e1
Force-Split $s2 $e2

# ---------------------------------------------------------------------
# Edit 3: " It can be seen on the example of modifying an event. The
#          changes are entered in a table in the view part of the
#          system (" + "Michał" + " will develop that in the
#          implementation part) and thanks to association, the changes
#          are also made on the particular event."
#      -> " It can be " / "explained " / "on the example of modifying
#          an event. The changes are entered in the view " / "and are
#          expected to be implemented in the " + "eventList" + " as well."
# ---------------------------------------------------------------------
$full = $d.Content.Text
$oldLead = " It can be seen on the example of modifying an event. The changes are entered in a table in the view part of the system ("
$anchor3 = $oldLead + "Michał"
$idx3 = $full.IndexOf($anchor3)
$r = $d.Range($idx3, $idx3 + $oldLead.Length)
$newLead = " It can be explained on the example of modifying an event. The changes are entered in the view and are expected to be implemented in the "
$r.Text = $newLead

$q1 = " It can be "
$q2 = "explained "
$q3 = "on the example of modifying an event. The changes are entered in the view "
$q4 = "and are expected to be implemented in the "

$qs1 = $idx3;  $qe1 = $qs1 + $q1.Length
$qs2 = $qe1;   $qe2 = $qs2 + $q2.Length
$qs3 = $qe2;   $qe3 = $qs3 + $q3.Length
$qs4 = $qe3;   $qe4 = $qs4 + $q4.Length

Force-Split $qs1 $qe1
Force-Split $qs2 $qe2
Force-Split $qs3 $qe3
Force-Split $qs4 $qe4

# "Michał" -> "eventList" (stays inside the existing proofErr run)
$full = $d.Content.Text
$idxMichal = $full.IndexOf("Michał will develop")
$rm = $d.Range($idxMichal, $idxMichal + "Michał".Length)
$rm.Text = "eventList"

# " will develop that in the implementation part) and thanks to
#   association, the changes are also made on the particular event."
# -> " as well."
$full = $d.Content.Text
$oldTail = " will develop that in the implementation part) and thanks to association, the changes are also made on the particular event."
$anchorTail = "eventList" + $oldTail
$idxTail = $full.IndexOf($anchorTail)
$tailStart = $idxTail + "eventList".Length
$rt = $d.Range($tailStart, $tailStart + $oldTail.Length)
$newTail = " as well."
$rt.Text = $newTail
$tailEnd = $tailStart + $newTail.Length
Force-Split $tailStart $tailEnd

# ---------------------------------------------------------------------
# Edit 4: split "...with more fields and methods." into
#         "...with more fiel" + "ds and methods." and move the
#         "_GoBack" bookmark to sit between the two new runs.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$anchor4 = "with more fields and methods."
$idx4 = $full.IndexOf($anchor4)
$splitPos = $idx4 + "with more fiel".Length
Force-SplitAt $splitPos

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
